# Apply the BOM.xlsx edits described by the commit:
# "Broke the BOM into individual files for Bob and Geoff to order stuff"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# --- Reassign some "Raw Materials" / "Mechanical" line items from Mike to Geoff ---
$ws.Range("C35").Value2 = "Geoff"
$ws.Range("C46").Value2 = "Geoff"
$ws.Range("C47").Value2 = "Geoff"
$ws.Range("C48").Value2 = "Geoff"
$ws.Range("C53").Value2 = "Geoff"

# --- Steel cut quantities: 4 cuts added for the 2x2 tubing rows ---
$ws.Range("O38").Value2 = 4
$ws.Range("O39").Value2 = 4
$ws.Range("O40").Value2 = 4
$ws.Range("O41").Value2 = 4

# --- The old (wrong) total-feet figure for row 43 is cleared ---
$ws.Range("N43").ClearContents()

# --- New "cut list" helper table in K53:P63 ---
$ws.Range("L53").Value2 = "10 foot lengths"
$ws.Range("N53").Value2 = "2x2 total"
$ws.Range("O53").Value2 = "3x3 total"
$ws.Range("P53").Value2 = "3x6 total"

$ws.Range("K54").Value2 = "cut1"
$ws.Range("L54").Formula = "=(M37+M38+M39)/12"
$ws.Range("M54").Value2 = "2x2"
$ws.Range("N54").Formula = "=L54+L55+L56"

$ws.Range("K55").Value2 = "cut2"
$ws.Range("L55").Formula = "=(M40+M41+M42)/12"
$ws.Range("M55").Value2 = "2x2"

$ws.Range("K56").Value2 = "cut3"
$ws.Range("L56").Formula = "=(M43+M44+M45+M46+M47+M48+M49+M50)/12"
$ws.Range("M56").Value2 = "2x2"

$ws.Range("K57").Value2 = "cut4"
$ws.Range("L57").Formula = "=(N37+N39)/12"
$ws.Range("M57").Value2 = "3x3"

$ws.Range("K58").Value2 = "cut5"
$ws.Range("L58").Formula = "=(N38+N40)/12"
$ws.Range("M58").Value2 = "3x3"

$ws.Range("K59").Value2 = "cut6"
$ws.Range("L59").Formula = "=(N41)/12"
$ws.Range("M59").Value2 = "3x3"

$ws.Range("K60").Value2 = "cut7"
$ws.Range("L60").Formula = "=N42/12"
$ws.Range("M60").Value2 = "3x3"

$ws.Range("K61").Value2 = "cut8"
$ws.Range("K62").Value2 = "cut9"
$ws.Range("K63").Value2 = "cut10"

$ws.Range("L57:L63").NumberFormat = "0"

# --- sheet view tweaks ---
$ws.Activate()
$ws.Range("D37").Select()
$excel.ActiveWindow.Zoom = 100

$wb.Application.CalculateFull()
